$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.532.60'
$ws.Range('D2').NumberFormat = 'General'
$ws.Range('E2').Value = '  -1.19%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.423.96'
$ws.Range('D3').NumberFormat = 'General'
$ws.Range('E3').Value = '  -2.14%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').NumberFormat = 'General'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '577.45'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  -2.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '128.85'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  -3.61%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('E8').Value = '  -1.50%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.53'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('E9').Value = '  +2.37%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.122'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E10').Value = '  -1.58%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.380'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('E11').Value = '  -1.24%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.008.68'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('E12').Value = '  -2.17%  '
$ws.Range('E14').Value = '  -3.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.427.71'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('E15').Value = '  -2.14%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.580.37'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '25.09'
$ws.Range('D17').NumberFormat = 'General'
$ws.Range('E17').Value = '  -2.14%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '9.81'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('E18').Value = '  -0.42%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '5.63'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('E19').Value = '  -2.06%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.24'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  -2.08%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '383.20'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').Value = '  -2.48%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.561'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  -1.95%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.562.79'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  -2.13%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '73.76'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  -1.16%  '
$ws.Range('E25').Value = '  +0.27%  '
$ws.Range('E26').Value = '  -5.23%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.02'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').Value = '  +2.34%  '
$ws.Range('E28').Value = '  -3.50%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.00'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('E29').Value = '  -4.87%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.88'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').Value = '  -4.01%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.41'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('E32').Value = '  -4.28%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.455.94'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E33').Value = '  -1.82%  '
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '22.63'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').Value = '  -3.44%  '
$ws.Range('E36').Value = '  +0.46%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.72'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  -2.28%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '164.17'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  -1.75%  '
$ws.Range('E39').Value = '  -2.90%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0762'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  -2.32%  '
$ws.Range('E41').Value = '  -3.09%  '
$ws.Range('E42').Value = '  +0.01%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '41.30'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').Value = '  -1.10%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.30'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('E44').Value = '  -2.10%  '
$ws.Range('E45').Value = '  -4.28%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '23.22'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').Value = '  -7.71%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.09'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('E47').Value = '  -6.36%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '6.69'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  -0.81%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.879'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('E49').Value = '  -1.30%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.267.80'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').Value = '  -4.56%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0250'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('E51').Value = '  -2.89%  '
